$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (row 15, column A has the bold/
# bordered/centered style) down onto the new rows so the new "A" cells inherit the
# same style (s="1") used by the rest of column A.
$ws.Range("A15").Copy()
$ws.Range("A16:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 16
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 15
$ws.Cells.Item(16, 3).Value = 0.19
$ws.Cells.Item(16, 4).Value = 2.8
$ws.Cells.Item(16, 5).Value = -0.22
$ws.Cells.Item(16, 6).Value = 3.39
$ws.Cells.Item(16, 7).Value = -0.08620689655172414
$ws.Cells.Item(16, 8).Value = 0.008474576271186441

# Row 17
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 16
$ws.Cells.Item(17, 3).Value = 0.44
$ws.Cells.Item(17, 4).Value = 3.23
$ws.Cells.Item(17, 5).Value = -0.24
$ws.Cells.Item(17, 6).Value = 3.71
$ws.Cells.Item(17, 7).Value = -0.1185185185185185
$ws.Cells.Item(17, 8).Value = 0.1467889908256881

# Row 18
$ws.Cells.Item(18, 1).Value = 16
$ws.Cells.Item(18, 2).Value = 17
$ws.Cells.Item(18, 3).Value = 0.23
$ws.Cells.Item(18, 4).Value = 2.26
$ws.Cells.Item(18, 5).Value = -0.14
$ws.Cells.Item(18, 6).Value = 2.49
$ws.Cells.Item(18, 7).Value = -0.1129032258064516
$ws.Cells.Item(18, 8).Value = 0

# Row 19
$ws.Cells.Item(19, 1).Value = 17
$ws.Cells.Item(19, 2).Value = 18
$ws.Cells.Item(19, 3).Value = 0.14
$ws.Cells.Item(19, 4).Value = 1.85
$ws.Cells.Item(19, 5).Value = 0.15
$ws.Cells.Item(19, 6).Value = 1.21
$ws.Cells.Item(19, 7).Value = -0.140625
$ws.Cells.Item(19, 8).Value = -0.1810344827586207
